$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D7").Value = 0.005
$ws.Range("E2:E7").Value = 0.06646196723119815
